# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (right after "总计") holding the
# quarterly fund-holding breakdown, and prepends a matching summary row
# to the "总计" sheet so it stays in sync with the new quarter.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $val) {
    # Force a literal/text value even when it "looks like" a number
    # (e.g. "005702", "0.20") so Excel doesn't silently coerce it and
    # drop the leading/trailing zeros. ClearFormats() afterwards removes
    # the NumberFormat="@" override we used to pin the type, so the cell
    # ends up with no explicit style -- matching the rest of the sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$q3.Range("A2").Value = 0
Set-TextValue $q3.Range("B2") "005702"
$q3.Range("C2").Value = "恒生前海港股通高股息低波动指数"
Set-TextValue $q3.Range("D2") "0.20"
Set-TextValue $q3.Range("E2") "94.22"
Set-TextValue $q3.Range("F2") "2.36"
Set-TextValue $q3.Range("G2") "0.0047"
$q3.Range("H2").Value = 9

# Re-apply the shared "bold + thin border + centered" header/index style
# (style index 2 in the original workbook) by copying it from an
# existing sheet that already carries it, rather than recreating it
# through Font/Border properties (which resolve to a different, richer
# default font in this host).
$headerStyleSrc = $wb.Worksheets.Item("2022-Q1").Range("B1")
$headerStyleSrc.Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$indexStyleSrc = $wb.Worksheets.Item("2022-Q1").Range("A2")
$indexStyleSrc.Copy()
$q3.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Prepend a "2022-Q3" row to the "总计" summary sheet
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0
# The row-insert drags a stray style onto the newly-blank B2:D2 cells;
# these should stay unstyled like the other data cells in the column.
$totalSheet.Range("B2:D2").ClearFormats()

$indexStyleSrc2 = $totalSheet.Range("A3")
$indexStyleSrc2.Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("A2").Value = 0

# Renumber the 0-based index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("A8").Value = 6
